# Updated cryptos list on Fri Jan 19 03:50:44 UTC 2024 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for every coin row,
# and swaps the ARBITRUM/Celestia rows (37/38) back to their scraped order.
#
# Price cells that look numeric ("0.998", "310.10", ...) are written with
# NumberFormat "@" (Text) so Excel's COM layer doesn't silently convert
# them to real numbers, then the style is reset to "Normal" so no stray
# cell-style reference is left behind - matching the plain text cells in
# the source file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.942.01"
$ws.Range("E2").Value = "  -3.90%  "
$ws.Range("D3").Value = "2.454.09"
$ws.Range("E3").Value = "  -3.00%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.22%  "
$ws.Range("E7").Value = "  -2.95%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.493"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.06"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0773"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.68%  "
$ws.Range("E12").Value = "  -0.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.61%  "
$ws.Range("D14").Value = "2.831.78"
$ws.Range("E14").Value = "  -2.94%  "
$ws.Range("D15").Value = "2.474.17"
$ws.Range("E15").Value = "  -2.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.72"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.775"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.07%  "
$ws.Range("D18").Value = "40.924.42"
$ws.Range("E18").Value = "  -3.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.68%  "
$ws.Range("D20").Value = "0.0₃0911"
$ws.Range("E20").Value = "  -3.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -9.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.49%  "
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.51"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.58"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "151.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.44"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.15%  "
$ws.Range("E33").Value = "  -5.79%  "
$ws.Range("E34").Value = "  -4.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0734"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.95%  "
$ws.Range("E36").Value = "  -5.26%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.84"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.06%  "
$ws.Range("B38").Value = "Celestia"
$ws.Range("C38").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.62"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.77%  "
$ws.Range("E39").Value = "  -3.56%  "
$ws.Range("E40").Value = "  -8.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.14"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.79%  "
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -12.04%  "
$ws.Range("D44").Value = "1.959.59"
$ws.Range("E44").Value = "  -2.40%  "
$ws.Range("E45").Value = "  -5.82%  "
$ws.Range("E46").Value = "  -8.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "69.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "96.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.176"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.14%  "
